# Implemented Ca40 input file
$wb = $excel.ActiveWorkbook

# --- Sheet1: tweak the epsilon parameter (E2) and the initial guess (G2) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("E2").Formula = "=0.04"
$ws1.Range("G2").Value = 300

# --- Sheet2: refine the SQRT() inputs for C2 / D2, move the selection ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("C2").Formula = "=SQRT(109.6264)"
$ws2.Range("D2").Formula = "=SQRT(190.4306)"
$ws2.Activate()
$ws2.Range("B5").Select()

# --- Sheet4 (the new Ca40 input) becomes the active sheet/tab ---
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Activate()
